$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 52-101: column A = medicine name (Chinese), column B = "<n>.jpg" filename,
# continuing the existing table pattern (header in row 1, data from row 2).
$names = @(
    "桑白皮",
    "枇杷葉",
    "茯苓",
    "白朮",
    "車前子",
    "木通",
    "澤瀉",
    "厚朴",
    "蒼朮",
    "砂仁",
    "威靈仙",
    "桑寄生",
    "獨活",
    "山楂",
    "麥芽",
    "大黃",
    "蘆薈",
    "火麻仁",
    "丹參",
    "桃仁",
    "紅花",
    "延胡索",
    "川芎",
    "薑黃",
    "益母草",
    "牛膝",
    "水蛭",
    "白及",
    "艾葉",
    "側柏葉",
    "三七",
    "金銀花",
    "連翹",
    "蒲公英",
    "敗醬草",
    "射干",
    "槴子",
    "夏枯草",
    "黃連",
    "黃岑",
    "黃柏",
    "龍膽",
    "苦參",
    "玄參",
    "牡丹皮",
    "紫草",
    "赤芍",
    "青蒿",
    "蛇床子",
    "檳榔"
)
$files = @(
    "51.jpg",
    "52.jpg",
    "53.jpg",
    "54.jpg",
    "55.jpg",
    "56.jpg",
    "57.jpg",
    "58.jpg",
    "59.jpg",
    "60.jpg",
    "61.jpg",
    "62.jpg",
    "63.jpg",
    "64.jpg",
    "65.jpg",
    "66.jpg",
    "67.jpg",
    "68.jpg",
    "69.jpg",
    "70.jpg",
    "71.jpg",
    "72.jpg",
    "73.jpg",
    "74.jpg",
    "75.jpg",
    "76.jpg",
    "77.jpg",
    "78.jpg",
    "79.jpg",
    "80.jpg",
    "81.jpg",
    "82.jpg",
    "83.jpg",
    "84.jpg",
    "85.jpg",
    "86.jpg",
    "87.jpg",
    "88.jpg",
    "89.jpg",
    "90.jpg",
    "91.jpg",
    "92.jpg",
    "93.jpg",
    "94.jpg",
    "95.jpg",
    "96.jpg",
    "97.jpg",
    "98.jpg",
    "99.jpg",
    "100.jpg"
)

$startRow = 52
# Template rows whose existing cell formatting we reuse so no duplicate styles are created:
#   column A formatting template -> row 32 (style used by A32:A51)
#   column B formatting template -> row 2  (style used by B2:B51)
$aFormatTemplate = $ws.Cells.Item(32, 1)
$bFormatTemplate = $ws.Cells.Item(2, 2)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $startRow + $i

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $names[$i]
    $aFormatTemplate.Copy() | Out-Null
    $aCell.PasteSpecial(-4122) | Out-Null

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = $files[$i]
    $bFormatTemplate.Copy() | Out-Null
    $bCell.PasteSpecial(-4122) | Out-Null

    $ws.Rows.Item($r).RowHeight = 19.95
}

# Row 99 (source string "青蒿") carries a distinct font (細明體, 16pt, black, family 3)
# instead of the usual column-A template font.
$specialCell = $ws.Cells.Item(99, 1)
$specialCell.Font.Name = "細明體"
$specialCell.Font.Size = 16
$specialCell.Font.Color = 0
$specialCell.Font.Family = 3

# Restore the clipboard / selection, and move the view roughly back to where the
# author left it (bottom of the newly extended table).
$excel.CutCopyMode = $false
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A102").Select() | Out-Null
